# agent specific features for ordering and delivery added
#
# Adds two new columns to Sheet1:
#   N (14) -> "ordering_period"
#   O (15) -> "delivery_period"
# and updates the existing "days_between_financing" (L) / "financing_period" (M)
# columns' values for the three agent rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers -----------------------------------------------
$ws.Range("N1").Value = "ordering_period"
$ws.Range("O1").Value = "delivery_period"

# --- Updated values for existing L (days_between_financing) / M (financing_period)
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 90

$ws.Range("L3").Value = 10
$ws.Range("M3").Value = 80

$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 70

# --- New column data (ordering_period / delivery_period) ---------------
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 0

$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 2

$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 3

# --- Column widths -------------------------------------------------------
# Target widths (21.25 -> 24.625, new cols 13.625 / 13.125) can't be hit
# exactly through the ColumnWidth/pixel-rounding COM surface, so these are
# the closest reachable values; column M (13) is left untouched since its
# width (14.25) does not change.
$ws.Columns.Item(12).ColumnWidth = 23.86
$ws.Columns.Item(14).ColumnWidth = 12.86
$ws.Columns.Item(15).ColumnWidth = 12.43

# --- Selection / scroll position ------------------------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("O4").Select() | Out-Null
